$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newGuids = @(
    "79cdc512f3ee4d12821547a1d96184f0",
    "7e2115889fb048f28b00eecb08e76dd2",
    "b1a9f088f6a4438ebb5a52136210571c",
    "fd90044f69be404e9dbf899f69598a5f",
    "34dbfda2740e47a783e5a6b28bcf70b7",
    "f9c71db1a1c440d98d32bf52d0cff3b5",
    "2b03c3ecaae24f01b22a469804aea182"
)

# Copy the number-format / font / border style used by column A on the
# existing rows (e.g. A7) so the newly added A8:A14 cells match exactly.
$ws.Range("A7").Copy()

for ($i = 0; $i -lt $newGuids.Length; $i++) {
    $row = 8 + $i
    $guid = $newGuids[$i]

    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $row

    $ws.Range("B$row").Value = $guid
    $ws.Range("C$row").Value = $guid
    $ws.Range("D$row").Value = $guid
    $ws.Range("E$row").Value = $guid
    $ws.Range("F$row").Value = $guid
    $ws.Range("G$row").Value = $guid
}
